$wb = $excel.ActiveWorkbook

$names = @(
    "summ11916577",
    "summ00587463",
    "summ01008681",
    "summ00987669",
    "summ59611014",
    "summ57973025",
    "summ55478695",
    "summ56904950",
    "summ58001593"
)

for ($i = 1; $i -le $names.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $names[$i - 1]
}
